# Add data for 2022-06-27: update the "through" date references and
# refresh the June (row 7) and Total (row 14) figures in the "2022"
# column (I).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-06-19"

# Update the column header text (shared string) for the 2022 total column.
$ws.Range("I1").Value = "2022 (through 06-19)"

# Update the June figure for 2022.
$ws.Range("I7").Value = 98

# Update the yearly total figure for 2022.
$ws.Range("I14").Value = 761
